$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.634.03"
$ws.Range("E2").Value = "  +5.43%  "
$ws.Range("D3").Value = "3.181.43"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'400.70"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").Value = "'109.95"
$ws.Range("E6").Value = "  +6.03%  "
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("D10").Value = "'39.20"
$ws.Range("E10").Value = "  +5.19%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "'0.0886"
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").Value = "3.680.66"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "'19.09"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "'8.05"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("E16").Value = "  +7.52%  "
$ws.Range("D17").Value = "3.184.29"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").Value = "'10.56"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "54.477.77"
$ws.Range("E19").Value = "  +4.91%  "
$ws.Range("D20").Value = "'3.30"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0000101"
$ws.Range("E21").Value = "  +4.03%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'12.91"
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("D23").Value = "'72.22"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'275.77"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").Value = "'3.26"
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("D26").Value = "'8.02"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").Value = "'27.85"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("E28").Value = "  +5.46%  "
$ws.Range("D29").Value = "'0.170"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").Value = "'11.05"
$ws.Range("E32").Value = "  +6.62%  "
$ws.Range("D33").Value = "'0.0515"
$ws.Range("E33").Value = "  +14.14%  "
$ws.Range("D34").Value = "'36.69"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("D37").Value = "'3.63"
$ws.Range("E37").Value = "  +6.28%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  +10.64%  "
$ws.Range("D40").Value = "'4.08"
$ws.Range("E40").Value = "  +10.18%  "
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").Value = "'0.293"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").Value = "'17.22"
$ws.Range("E43").Value = "  +1.84%  "
$ws.Range("D44").Value = "'131.08"
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "'22.11"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "'2.49"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "2.095.72"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("E50").Value = "  +6.45%  "
$ws.Range("E51").Value = "  +13.74%  "
